# Scheduled data refresh: update currentAveragePrice-derived Leve profit figures
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1058.8611
$ws.Range("J17").Value = 1107.7097
$ws.Range("L17").Value = 3323.1291
$ws.Range("N17").Value = -3659.1291
$ws.Range("H41").Value = 639.8333
$ws.Range("I41").Value = 284.75
$ws.Range("K41").Value = 284.75
$ws.Range("M41").Value = 155.25
$ws.Range("H74").Value = 8398.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8398.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 8398.4
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -10270.4
$ws.Range("H77").Value = 8398.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8398.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 41992
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -51352
$ws.Range("H131").Value = 3025.9546
$ws.Range("J131").Value = 6424.5713
$ws.Range("L131").Value = 19273.7139
$ws.Range("N131").Value = -29353.7139
$ws.Range("H137").Value = 9031.558999999999
$ws.Range("I137").Value = 15374.6875
$ws.Range("J137").Value = 3393.2222
$ws.Range("K137").Value = 46124.0625
$ws.Range("L137").Value = 10179.6666
$ws.Range("M137").Value = -43574.0625
$ws.Range("N137").Value = -15279.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7663.364
$ws.Range("I45").Value = 7659.6
$ws.Range("K45").Value = 7659.6
$ws.Range("M45").Value = -7282.6
$ws.Range("H122").Value = 2004735.9
$ws.Range("I122").Value = 4606.1333
$ws.Range("K122").Value = 13818.3999
$ws.Range("M122").Value = -11368.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 8512.0625
$ws.Range("I11").Value = 9869.5
$ws.Range("J11").Value = 6249.6665
$ws.Range("K11").Value = 9869.5
$ws.Range("L11").Value = 6249.6665
$ws.Range("M11").Value = -9729.5
$ws.Range("N11").Value = -6529.6665
$ws.Range("H20").Value = 4942.3335
$ws.Range("I20").Value = 1930.8
$ws.Range("K20").Value = 1930.8
$ws.Range("M20").Value = -1683.8
$ws.Range("H82").Value = 59912
$ws.Range("J82").Value = 107499
$ws.Range("L82").Value = 107499
$ws.Range("N82").Value = -108265
$ws.Range("H85").Value = 59912
$ws.Range("J85").Value = 107499
$ws.Range("L85").Value = 107499
$ws.Range("N85").Value = -110151
$ws.Range("H134").Value = 10963
$ws.Range("I134").Value = 11856.84
$ws.Range("J134").Value = 6493.8
$ws.Range("K134").Value = 35570.52
$ws.Range("L134").Value = 19481.4
$ws.Range("M134").Value = -33035.52
$ws.Range("N134").Value = -24551.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 100000
$ws.Range("J9").Value = 100000
$ws.Range("L9").Value = 100000
$ws.Range("N9").Value = -100336
$ws.Range("H22").Value = 673.6
$ws.Range("I22").Value = 678.7778
$ws.Range("K22").Value = 678.7778
$ws.Range("M22").Value = -328.7778
$ws.Range("H58").Value = 3210.75
$ws.Range("I58").Value = 2039.8
$ws.Range("J58").Value = 4047.1428
$ws.Range("K58").Value = 2039.8
$ws.Range("L58").Value = 4047.1428
$ws.Range("M58").Value = -1836.8
$ws.Range("N58").Value = -4453.1428
$ws.Range("H136").Value = 3210.75
$ws.Range("I136").Value = 2039.8
$ws.Range("J136").Value = 4047.1428
$ws.Range("K136").Value = 6119.4
$ws.Range("L136").Value = 12141.4284
$ws.Range("M136").Value = -3569.4
$ws.Range("N136").Value = -17241.4284
$ws.Range("H138").Value = 146428.14
$ws.Range("J138").Value = 146428.14
$ws.Range("L138").Value = 146428.14
$ws.Range("N138").Value = -156708.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2348.8
$ws.Range("I70").Value = 1561
$ws.Range("K70").Value = 4683
$ws.Range("M70").Value = -4368
$ws.Range("H73").Value = 2348.8
$ws.Range("I73").Value = 1561
$ws.Range("K73").Value = 4683
$ws.Range("M73").Value = -3591
$ws.Range("H97").Value = 45589.473
$ws.Range("J97").Value = 40833.5
$ws.Range("L97").Value = 122500.5
$ws.Range("N97").Value = -123492.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 10000
$ws.Range("M5").Value = -9888
$ws.Range("H57").Value = 33997.8
$ws.Range("J57").Value = 39994.5
$ws.Range("L57").Value = 39994.5
$ws.Range("N57").Value = -41634.5
$ws.Range("H70").Value = 6726.6875
$ws.Range("I70").Value = 7432.2856
$ws.Range("J70").Value = 6177.8887
$ws.Range("K70").Value = 7432.2856
$ws.Range("L70").Value = 6177.8887
$ws.Range("M70").Value = -7162.2856
$ws.Range("N70").Value = -6717.8887
$ws.Range("H73").Value = 6726.6875
$ws.Range("I73").Value = 7432.2856
$ws.Range("J73").Value = 6177.8887
$ws.Range("K73").Value = 7432.2856
$ws.Range("L73").Value = 6177.8887
$ws.Range("M73").Value = -6496.2856
$ws.Range("N73").Value = -8049.8887
$ws.Range("H141").Value = 99997
$ws.Range("J141").Value = 99997
$ws.Range("L141").Value = 99997
$ws.Range("N141").Value = -110357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 606.4
$ws.Range("I9").Value = 570.5
$ws.Range("K9").Value = 570.5
$ws.Range("M9").Value = -346.5
$ws.Range("H22").Value = 9319.706
$ws.Range("I22").Value = 9975.532999999999
$ws.Range("J22").Value = 4401
$ws.Range("K22").Value = 9975.532999999999
$ws.Range("L22").Value = 4401
$ws.Range("M22").Value = -9680.532999999999
$ws.Range("N22").Value = -4991
$ws.Range("H27").Value = 9319.706
$ws.Range("I27").Value = 9975.532999999999
$ws.Range("J27").Value = 4401
$ws.Range("K27").Value = 9975.532999999999
$ws.Range("L27").Value = 4401
$ws.Range("M27").Value = -9868.532999999999
$ws.Range("N27").Value = -4615
$ws.Range("H100").Value = 2314.7273
$ws.Range("I100").Value = 1923.2858
$ws.Range("J100").Value = 2999.75
$ws.Range("K100").Value = 1923.2858
$ws.Range("L100").Value = 2999.75
$ws.Range("M100").Value = -1382.2858
$ws.Range("N100").Value = -4081.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 19375
$ws.Range("I47").Value = 18000
$ws.Range("K47").Value = 18000
$ws.Range("M47").Value = -17428
$ws.Range("H126").Value = 18530.115
$ws.Range("I126").Value = 25217.47
$ws.Range("J126").Value = 5898.4443
$ws.Range("K126").Value = 75652.41
$ws.Range("L126").Value = 17695.3329
$ws.Range("M126").Value = -73182.41
$ws.Range("N126").Value = -22635.3329

Write-Host "Updated Leve profit values across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
